$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New column E (OCV) data for rows 3-24
$ws.Range("E3").Value = 2.98
$ws.Range("E4").Value = 3.177
$ws.Range("E5").Value = 3.227
$ws.Range("E6").Value = 3.25
$ws.Range("E7").Value = 3.273
$ws.Range("E8").Value = 3.29
$ws.Range("E9").Value = 3.301
$ws.Range("E10").Value = 3.308
$ws.Range("E11").Value = 3.311
$ws.Range("E12").Value = 3.312
$ws.Range("E13").Value = 3.314
$ws.Range("E14").Value = 3.316
$ws.Range("E15").Value = 3.317
$ws.Range("E16").Value = 3.319
$ws.Range("E17").Value = 3.321
$ws.Range("E18").Value = 3.324
$ws.Range("E19").Value = 3.328
$ws.Range("E20").Value = 3.338
$ws.Range("E21").Value = 3.345
$ws.Range("E22").Value = 3.347
$ws.Range("E23").Value = 3.349
$ws.Range("E24").Value = 3.352

# Update selection to E24
$ws.Range("E24").Select()

# Update window size
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Width = 21600
$wb.Windows.Item(1).Height = 11620
